$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f76a74000c6e6bf7bba5f470370dfbe35446590/e2e/c09e2566-92f6-4af4-9a07-8febdfd04abe.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a1c4cc2eb721fb97b6b5cc88050e49a1f9526a4/e2e/c09e2566-92f6-4af4-9a07-8febdfd04abe.md."
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a1c4cc2eb721fb97b6b5cc88050e49a1f9526a4/e2e/c09e2566-92f6-4af4-9a07-8febdfd04abe.md"
$displayName = "c09e2566-92f6-4af4-9a07-8febdfd04abe.md"

# ---- zh-cn sheet: row 7 (c09e2566-...) gets a handback recorded for it ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = $wsZh.Range("G7").Value2
$wsZh.Range("K7").Value = "2016-09-06 11:13:51"
$wsZh.Range("P7").Value = $errorMessage

$wsZh.Range("I7").Value = $displayName
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackUrl, "", "", $displayName)

# ---- de-de sheet: row 7 (c09e2566-...) gets a handback recorded for it ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = $wsDe.Range("G7").Value2
$wsDe.Range("K7").Value = "2016-09-06 11:13:59"
$wsDe.Range("P7").Value = $errorMessage

$wsDe.Range("I7").Value = $displayName
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackUrl, "", "", $displayName)
